$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new MultiWOZ dataset row (row 12) mirroring the structure of the
# existing rows (columns: Name, Introduction, Multi/Single Turn, Task,
# Task Detail, Public Accessible, Links, Size & Stats, Included Label,
# Missing Label).
$ws.Range("A12").Value = "MultiWOZ"
$ws.Range("B12").Value = "EMNLP 2018 best paper, not release yet."
$ws.Range("C12").Value = "N/A"
$ws.Range("D12").Value = "N/A"
$ws.Range("E12").Value = "N/A"
$ws.Range("F12").Value = "N/A"
$ws.Range("G12").Value = "Paper:`nMultiWOZ - A Large-Scale Multi-Domain Wizard-of-Oz Dataset for Task-Oriented Dialogue Modelling（Not Released Yet）"
$ws.Range("H12").Value = "N/A"
$ws.Range("I12").Value = "N/A"
$ws.Range("J12").Value = "N/A"

# Match the row height of the preceding data rows (57pt, same as rows 10-11).
$ws.Rows.Item(12).RowHeight = 57
